$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1): I0, IF - reuse H1's header style (bold, bordered,
# centered) by copying its format rather than rebuilding it from scratch.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells (row 2): plain numbers, default style
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 7
